# Update Fecha (D), Volumen (J), Precio mínimo (K), Precio máximo (L),
# Precio promedio ponderado (M) and Precio $/Kg (P) for the weekly data rows.
# Values got reshuffled between rows (same data set, different row order).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2  = @{ D = 44221; J = 50; K = 2500; L = 2500; M = 2500; P = 833 }
    4  = @{ D = 44340; J = 54; K = 3000; L = 3000; M = 3000; P = 1000 }
    5  = @{ D = 44222; J = 45; K = 3000; L = 3000; M = 3000; P = 1000 }
    6  = @{ D = 44260; J = 60; K = 3500; L = 3500; M = 3500; P = 1167 }
    8  = @{ D = 44224; J = 67; K = 3000; L = 3000; M = 3000; P = 1000 }
    9  = @{ D = 44291; J = 45; K = 3000; L = 3000; M = 3000; P = 1000 }
    10 = @{ D = 44292; J = 40; K = 3000; L = 3000; M = 3000; P = 1000 }
    11 = @{ D = 44165; J = 68; K = 3000; L = 3000; M = 3000; P = 1000 }
    12 = @{ D = 44243; J = 45; K = 3000; L = 3000; M = 3000; P = 1000 }
    13 = @{ D = 44242; J = 95; K = 2500; L = 3000; M = 2737; P = 912 }
    14 = @{ D = 44166; J = 45; K = 2500; L = 2500; M = 2500; P = 833 }
    15 = @{ D = 44223; J = 80; K = 2500; L = 3000; M = 2781; P = 927 }
    16 = @{ D = 44390; J = 50; K = 3000; L = 3000; M = 3000; P = 1000 }
    17 = @{ D = 44187; J = 65; K = 3000; L = 3000; M = 3000; P = 1000 }
    18 = @{ D = 44193; J = 70; K = 3000; L = 3000; M = 3000; P = 1000 }
    19 = @{ D = 44389; J = 81; K = 2800; L = 3000; M = 2889; P = 963 }
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 4).Value = $vals.D    # D: Fecha
    $ws.Cells.Item($row, 10).Value = $vals.J   # J: Volumen
    $ws.Cells.Item($row, 11).Value = $vals.K   # K: Precio minimo
    $ws.Cells.Item($row, 12).Value = $vals.L   # L: Precio maximo
    $ws.Cells.Item($row, 13).Value = $vals.M   # M: Precio promedio ponderado
    $ws.Cells.Item($row, 16).Value = $vals.P   # P: Precio $/Kg
}
